# Update column G ("K" = strikeouts) values for rows 2-35 on Sheet1.
# These values were regenerated from the source box-score data (K instead
# of the previous Strike# derived figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 1
    3  = 4
    4  = 4
    5  = 1
    6  = 4
    7  = 2
    8  = 7
    9  = 7
    10 = 2
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 3
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 3
    26 = 8
    27 = 3
    28 = 1
    29 = 2
    30 = 1
    31 = 0
    32 = 3
    33 = 4
    34 = 3
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
